$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 1100
$ws.Range("H23").Value = 1100
$ws.Range("H55").Value = 183.5
$ws.Range("I55").Value = 149.71428
$ws.Range("K55").Value = 149.71428
$ws.Range("M55").Value = 64.28572
$ws.Range("H61").Value = 6624.625
$ws.Range("I61").Value = 6624.625
$ws.Range("K61").Value = 19873.875
$ws.Range("M61").Value = -19701.875
$ws.Range("H76").Value = 5190.143
$ws.Range("J76").Value = 5996.5
$ws.Range("L76").Value = 5996.5
$ws.Range("N76").Value = -6626.5
$ws.Range("H79").Value = 5190.143
$ws.Range("J79").Value = 5996.5
$ws.Range("L79").Value = 5996.5
$ws.Range("N79").Value = -8180.5
$ws.Range("H80").Value = 1311.8928
$ws.Range("I80").Value = 1033.6666
$ws.Range("J80").Value = 1632.9231
$ws.Range("K80").Value = 3100.9998
$ws.Range("L80").Value = 4898.7693
$ws.Range("M80").Value = -2102.9998
$ws.Range("N80").Value = -6894.7693
$ws.Range("H83").Value = 1311.8928
$ws.Range("I83").Value = 1033.6666
$ws.Range("J83").Value = 1632.9231
$ws.Range("K83").Value = 9302.999400000001
$ws.Range("L83").Value = 14696.3079
$ws.Range("M83").Value = -4310.999400000001
$ws.Range("N83").Value = -24680.3079
$ws.Range("H88").Value = 3979.9
$ws.Range("J88").Value = 4519.8
$ws.Range("L88").Value = 4519.8
$ws.Range("N88").Value = -5331.8
$ws.Range("H91").Value = 3979.9
$ws.Range("J91").Value = 4519.8
$ws.Range("L91").Value = 4519.8
$ws.Range("N91").Value = -7327.8
$ws.Range("H98").Value = 5609364.5
$ws.Range("I98").Value = 6995889.5
$ws.Range("K98").Value = 6995889.5
$ws.Range("M98").Value = -6994391.5
$ws.Range("H107").Value = 33335766
$ws.Range("I107").Value = 22729320
$ws.Range("K107").Value = 22729320
$ws.Range("M107").Value = -22727400
$ws.Range("H112").Value = 4358044.5
$ws.Range("J112").Value = 4980321.5
$ws.Range("L112").Value = 14940964.5
$ws.Range("N112").Value = -14943180.5
$ws.Range("H121").Value = 6066.5
$ws.Range("J121").Value = 6066.5
$ws.Range("L121").Value = 18199.5
$ws.Range("N121").Value = -21693.5
$ws.Range("H122").Value = 5609364.5
$ws.Range("I122").Value = 6995889.5
$ws.Range("K122").Value = 20987668.5
$ws.Range("M122").Value = -20985218.5
$ws.Range("H127").Value = 250001680
$ws.Range("I127").Value = 142859060
$ws.Range("K127").Value = 428577180
$ws.Range("M127").Value = -428572220
$ws.Range("H131").Value = 18874.572
$ws.Range("I131").Value = 19020.416
$ws.Range("K131").Value = 57061.24800000001
$ws.Range("M131").Value = -52021.24800000001
$ws.Range("H132").Value = 3559.4211
$ws.Range("I132").Value = 3210.3794
$ws.Range("J132").Value = 4684.1113
$ws.Range("K132").Value = 9631.138199999999
$ws.Range("L132").Value = 14052.3339
$ws.Range("M132").Value = -7101.138199999999
$ws.Range("N132").Value = -19112.3339
$ws.Range("H137").Value = 40473.9
$ws.Range("I137").Value = 53694
$ws.Range("J137").Value = 9627
$ws.Range("K137").Value = 161082
$ws.Range("L137").Value = 28881
$ws.Range("M137").Value = -158532
$ws.Range("N137").Value = -33981
$ws.Range("H138").Value = 2088.6
$ws.Range("I138").Value = 867.6977000000001
$ws.Range("J138").Value = 3009.6316
$ws.Range("K138").Value = 2603.0931
$ws.Range("L138").Value = 9028.8948
$ws.Range("M138").Value = 2536.9069
$ws.Range("N138").Value = -19308.8948
$ws.Range("H141").Value = 2454.1177
$ws.Range("I141").Value = 2494.2144
$ws.Range("J141").Value = 2267
$ws.Range("K141").Value = 7482.6432
$ws.Range("L141").Value = 6801
$ws.Range("M141").Value = -2302.6432
$ws.Range("N141").Value = -17161
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2104.182
$ws.Range("I2").Value = 1931.7142
$ws.Range("K2").Value = 1931.7142
$ws.Range("M2").Value = -1818.7142
$ws.Range("H16").Value = 1595.75
$ws.Range("I16").Value = 1595.75
$ws.Range("K16").Value = 1595.75
$ws.Range("M16").Value = -1308.75
$ws.Range("H32").Value = 5015.0884
$ws.Range("I32").Value = 4125.082
$ws.Range("J32").Value = 12770.857
$ws.Range("K32").Value = 4125.082
$ws.Range("L32").Value = 12770.857
$ws.Range("M32").Value = -3838.082
$ws.Range("N32").Value = -13344.857
$ws.Range("H45").Value = 2796
$ws.Range("I45").Value = 2861.375
$ws.Range("K45").Value = 2861.375
$ws.Range("M45").Value = -2484.375
$ws.Range("H61").Value = 8951.923000000001
$ws.Range("I61").Value = 6851.7144
$ws.Range("J61").Value = 11402.167
$ws.Range("K61").Value = 6851.7144
$ws.Range("L61").Value = 11402.167
$ws.Range("M61").Value = -6639.7144
$ws.Range("N61").Value = -11826.167
$ws.Range("H74").Value = 111217.336
$ws.Range("I74").Value = 132724.8
$ws.Range("J74").Value = 3680
$ws.Range("K74").Value = 132724.8
$ws.Range("L74").Value = 3680
$ws.Range("M74").Value = -131850.8
$ws.Range("N74").Value = -5428
$ws.Range("H77").Value = 111217.336
$ws.Range("I77").Value = 132724.8
$ws.Range("J77").Value = 3680
$ws.Range("K77").Value = 663624
$ws.Range("L77").Value = 18400
$ws.Range("M77").Value = -659256
$ws.Range("N77").Value = -27136
$ws.Range("H88").Value = 2357.7
$ws.Range("J88").Value = 2397.8333
$ws.Range("L88").Value = 2397.8333
$ws.Range("N88").Value = -3209.8333
$ws.Range("H91").Value = 2357.7
$ws.Range("J91").Value = 2397.8333
$ws.Range("L91").Value = 2397.8333
$ws.Range("N91").Value = -5205.8333
$ws.Range("H97").Value = 639.86957
$ws.Range("I97").Value = 487.13635
$ws.Range("K97").Value = 487.13635
$ws.Range("M97").Value = 8.863650000000007
$ws.Range("H102").Value = 2332.1365
$ws.Range("I102").Value = 2285.35
$ws.Range("J102").Value = 2800
$ws.Range("K102").Value = 2285.35
$ws.Range("L102").Value = 2800
$ws.Range("M102").Value = -663.3499999999999
$ws.Range("N102").Value = -6044
$ws.Range("H110").Value = 4018.524
$ws.Range("I110").Value = 1389
$ws.Range("J110").Value = 15194
$ws.Range("K110").Value = 1389
$ws.Range("L110").Value = 15194
$ws.Range("M110").Value = 656
$ws.Range("N110").Value = -19284
$ws.Range("H116").Value = 2104.182
$ws.Range("I116").Value = 1931.7142
$ws.Range("K116").Value = 1931.7142
$ws.Range("M116").Value = 362.2858000000001
$ws.Range("H122").Value = 3143.4443
$ws.Range("I122").Value = 2771.6365
$ws.Range("K122").Value = 8314.9095
$ws.Range("M122").Value = -5864.9095
$ws.Range("H132").Value = 2393.1843
$ws.Range("I132").Value = 2480.9583
$ws.Range("K132").Value = 7442.874899999999
$ws.Range("M132").Value = -4912.874899999999
$ws.Range("H136").Value = 8951.923000000001
$ws.Range("I136").Value = 6851.7144
$ws.Range("J136").Value = 11402.167
$ws.Range("K136").Value = 20555.1432
$ws.Range("L136").Value = 34206.501
$ws.Range("M136").Value = -18005.1432
$ws.Range("N136").Value = -39306.501
$ws.Range("H140").Value = 51485.8
$ws.Range("J140").Value = 51485.8
$ws.Range("L140").Value = 51485.8
$ws.Range("N140").Value = -61845.8
$ws.Range("H141").Value = 51929
$ws.Range("J141").Value = 51929
$ws.Range("L141").Value = 51929
$ws.Range("N141").Value = -62289
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2104.182
$ws.Range("I3").Value = 1931.7142
$ws.Range("K3").Value = 1931.7142
$ws.Range("M3").Value = -1817.7142
$ws.Range("H22").Value = 990.6667
$ws.Range("I22").Value = 990.6667
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 990.6667
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -817.6667
$ws.Range("N22").ClearContents()
$ws.Range("H43").Value = 320000
$ws.Range("J43").Value = 320000
$ws.Range("L43").Value = 320000
$ws.Range("N43").Value = -320362
$ws.Range("H99").Value = 4768.72
$ws.Range("I99").Value = 5178.75
$ws.Range("J99").Value = 4575.7646
$ws.Range("K99").Value = 5178.75
$ws.Range("L99").Value = 4575.7646
$ws.Range("M99").Value = -3680.75
$ws.Range("N99").Value = -7571.7646
$ws.Range("H107").Value = 4030.4443
$ws.Range("I107").Value = 3909.25
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 3909.25
$ws.Range("L107").Value = 5000
$ws.Range("M107").Value = -1989.25
$ws.Range("N107").Value = -8840
$ws.Range("H115").Value = 35000
$ws.Range("J115").Value = 35000
$ws.Range("L115").Value = 35000
$ws.Range("N115").Value = -38134
$ws.Range("H134").Value = 2707
$ws.Range("I134").Value = 1573.091
$ws.Range("K134").Value = 4719.272999999999
$ws.Range("M134").Value = -2184.272999999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 296155.8
$ws.Range("I31").Value = 436178.25
$ws.Range("J31").Value = 3381.6365
$ws.Range("K31").Value = 436178.25
$ws.Range("L31").Value = 3381.6365
$ws.Range("M31").Value = -435883.25
$ws.Range("N31").Value = -3971.6365
$ws.Range("H34").Value = 296155.8
$ws.Range("I34").Value = 436178.25
$ws.Range("J34").Value = 3381.6365
$ws.Range("K34").Value = 436178.25
$ws.Range("L34").Value = 3381.6365
$ws.Range("M34").Value = -435976.25
$ws.Range("N34").Value = -3785.6365
$ws.Range("H58").Value = 2326.7827
$ws.Range("I58").Value = 2070.25
$ws.Range("K58").Value = 2070.25
$ws.Range("M58").Value = -1867.25
$ws.Range("H59").Value = 98000
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 98000
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 98000
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -100290
$ws.Range("H68").Value = 37396.668
$ws.Range("J68").Value = 37396.668
$ws.Range("L68").Value = 37396.668
$ws.Range("N68").Value = -38894.668
$ws.Range("H71").Value = 37396.668
$ws.Range("J71").Value = 37396.668
$ws.Range("L71").Value = 112190.004
$ws.Range("N71").Value = -119678.004
$ws.Range("H74").Value = 41056.8
$ws.Range("I74").Value = 32642.5
$ws.Range("J74").Value = 46666.332
$ws.Range("K74").Value = 32642.5
$ws.Range("L74").Value = 46666.332
$ws.Range("M74").Value = -31768.5
$ws.Range("N74").Value = -48414.332
$ws.Range("H77").Value = 41056.8
$ws.Range("I77").Value = 32642.5
$ws.Range("J77").Value = 46666.332
$ws.Range("K77").Value = 97927.5
$ws.Range("L77").Value = 139998.996
$ws.Range("M77").Value = -93559.5
$ws.Range("N77").Value = -148734.996
$ws.Range("H86").Value = 562173.6
$ws.Range("I86").Value = 775153.5600000001
$ws.Range("J86").Value = 8425.799999999999
$ws.Range("K86").Value = 775153.5600000001
$ws.Range("L86").Value = 8425.799999999999
$ws.Range("M86").Value = -774030.5600000001
$ws.Range("N86").Value = -10671.8
$ws.Range("H89").Value = 562173.6
$ws.Range("I89").Value = 775153.5600000001
$ws.Range("J89").Value = 8425.799999999999
$ws.Range("K89").Value = 3875767.8
$ws.Range("L89").Value = 42129
$ws.Range("M89").Value = -3870151.8
$ws.Range("N89").Value = -53361
$ws.Range("H105").Value = 4297.773
$ws.Range("I105").Value = 1550.1904
$ws.Range("K105").Value = 1550.1904
$ws.Range("M105").Value = 196.8096
$ws.Range("H107").Value = 5028.853
$ws.Range("I107").Value = 887.44446
$ws.Range("J107").Value = 6519.76
$ws.Range("K107").Value = 887.44446
$ws.Range("L107").Value = 6519.76
$ws.Range("M107").Value = 1032.55554
$ws.Range("N107").Value = -10359.76
$ws.Range("H122").Value = 5250.75
$ws.Range("I122").Value = 4973
$ws.Range("K122").Value = 14919
$ws.Range("M122").Value = -12469
$ws.Range("H134").Value = 6115.0356
$ws.Range("I134").Value = 6888.3184
$ws.Range("J134").Value = 3279.6667
$ws.Range("K134").Value = 20664.9552
$ws.Range("L134").Value = 9839.000100000001
$ws.Range("M134").Value = -18129.9552
$ws.Range("N134").Value = -14909.0001
$ws.Range("H136").Value = 2326.7827
$ws.Range("I136").Value = 2070.25
$ws.Range("K136").Value = 6210.75
$ws.Range("M136").Value = -3660.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 989.6
$ws.Range("I5").Value = 992
$ws.Range("J5").Value = 980
$ws.Range("K5").Value = 2976
$ws.Range("L5").Value = 2940
$ws.Range("M5").Value = -2864
$ws.Range("N5").Value = -3164
$ws.Range("H6").Value = 89
$ws.Range("I6").Value = 108.333336
$ws.Range("K6").Value = 325.000008
$ws.Range("M6").Value = -212.000008
$ws.Range("H17").Value = 3398.6316
$ws.Range("I17").Value = 2726.4546
$ws.Range("J17").Value = 4322.875
$ws.Range("K17").Value = 8179.3638
$ws.Range("L17").Value = 12968.625
$ws.Range("M17").Value = -8010.3638
$ws.Range("N17").Value = -13306.625
$ws.Range("H38").Value = 108.833336
$ws.Range("I38").Value = 29.5
$ws.Range("K38").Value = 88.5
$ws.Range("M38").Value = 258.5
$ws.Range("H44").Value = 2039.1818
$ws.Range("I44").Value = 681.1667
$ws.Range("J44").Value = 2548.4375
$ws.Range("K44").Value = 2043.5001
$ws.Range("L44").Value = 7645.3125
$ws.Range("M44").Value = -1645.5001
$ws.Range("N44").Value = -8441.3125
$ws.Range("H64").Value = 27781194
$ws.Range("I64").Value = 3399.5
$ws.Range("J64").Value = 83336780
$ws.Range("K64").Value = 10198.5
$ws.Range("L64").Value = 250010340
$ws.Range("M64").Value = -9928.5
$ws.Range("N64").Value = -250010880
$ws.Range("H67").Value = 27781194
$ws.Range("I67").Value = 3399.5
$ws.Range("J67").Value = 83336780
$ws.Range("K67").Value = 10198.5
$ws.Range("L67").Value = 250010340
$ws.Range("M67").Value = -9262.5
$ws.Range("N67").Value = -250012212
$ws.Range("H114").Value = 484
$ws.Range("I114").Value = 484
$ws.Range("K114").Value = 1452
$ws.Range("M114").Value = 1802
$ws.Range("H131").Value = 16130492
$ws.Range("I131").Value = 100000910
$ws.Range("J131").Value = 1564.5
$ws.Range("K131").Value = 300002730
$ws.Range("L131").Value = 4693.5
$ws.Range("M131").Value = -299997690
$ws.Range("N131").Value = -14773.5
$ws.Range("H134").Value = 7499.2173
$ws.Range("I134").Value = 7920.1055
$ws.Range("K134").Value = 23760.3165
$ws.Range("M134").Value = -18690.3165
$ws.Range("H135").Value = 989.6
$ws.Range("I135").Value = 992
$ws.Range("J135").Value = 980
$ws.Range("K135").Value = 8928
$ws.Range("L135").Value = 8820
$ws.Range("M135").Value = -6393
$ws.Range("N135").Value = -13890
$ws.Range("H140").Value = 27791772
$ws.Range("I140").Value = 55578212
$ws.Range("J140").Value = 5333.3335
$ws.Range("K140").Value = 166734636
$ws.Range("L140").Value = 16000.0005
$ws.Range("M140").Value = -166729456
$ws.Range("N140").Value = -26360.0005
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 18998
$ws.Range("J20").Value = 18998
$ws.Range("L20").Value = 18998
$ws.Range("N20").Value = -19488
$ws.Range("H24").Value = 13997
$ws.Range("J24").Value = 13997
$ws.Range("L24").Value = 13997
$ws.Range("N24").Value = -14343
$ws.Range("H26").Value = 40562.145
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 40562.145
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 40562.145
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -41122.145
$ws.Range("H50").Value = 40562.145
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 40562.145
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 40562.145
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -41558.145
$ws.Range("H70").Value = 5517.8887
$ws.Range("I70").Value = 5260.2705
$ws.Range("J70").Value = 6078.5884
$ws.Range("K70").Value = 5260.2705
$ws.Range("L70").Value = 6078.5884
$ws.Range("M70").Value = -4990.2705
$ws.Range("N70").Value = -6618.5884
$ws.Range("H73").Value = 5517.8887
$ws.Range("I73").Value = 5260.2705
$ws.Range("J73").Value = 6078.5884
$ws.Range("K73").Value = 5260.2705
$ws.Range("L73").Value = 6078.5884
$ws.Range("M73").Value = -4324.2705
$ws.Range("N73").Value = -7950.5884
$ws.Range("H102").Value = 49320.5
$ws.Range("I102").Value = 4051.85
$ws.Range("K102").Value = 4051.85
$ws.Range("M102").Value = -2429.85
$ws.Range("H113").Value = 8361.875
$ws.Range("J113").Value = 4413.7144
$ws.Range("L113").Value = 4413.7144
$ws.Range("N113").Value = -8753.714400000001
$ws.Range("H122").Value = 9144.861999999999
$ws.Range("I122").Value = 7911.222
$ws.Range("K122").Value = 23733.666
$ws.Range("M122").Value = -21283.666
$ws.Range("H132").Value = 24707.625
$ws.Range("I132").Value = 26617.727
$ws.Range("K132").Value = 79853.181
$ws.Range("M132").Value = -77323.181
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6014
$ws.Range("I7").Value = 6515.1924
$ws.Range("K7").Value = 6515.1924
$ws.Range("M7").Value = -6403.1924
$ws.Range("H43").Value = 12033.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 12033.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 12033.5
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -12419.5
$ws.Range("H46").Value = 1865.25
$ws.Range("J46").Value = 1735.3334
$ws.Range("L46").Value = 1735.3334
$ws.Range("N46").Value = -2111.3334
$ws.Range("H55").Value = 5647.4517
$ws.Range("I55").Value = 805
$ws.Range("J55").Value = 19569.5
$ws.Range("K55").Value = 805
$ws.Range("L55").Value = 19569.5
$ws.Range("M55").Value = -632
$ws.Range("N55").Value = -19915.5
$ws.Range("H93").Value = 50001804
$ws.Range("J93").Value = 200003120
$ws.Range("L93").Value = 200003120
$ws.Range("N93").Value = -200005616
$ws.Range("H122").Value = 3722.6667
$ws.Range("I122").Value = 3895.3333
$ws.Range("J122").Value = 3550
$ws.Range("K122").Value = 11685.9999
$ws.Range("L122").Value = 10650
$ws.Range("M122").Value = -9235.999899999999
$ws.Range("N122").Value = -15550
$ws.Range("H126").Value = 6014
$ws.Range("I126").Value = 6515.1924
$ws.Range("K126").Value = 19545.5772
$ws.Range("M126").Value = -17075.5772
$ws.Range("H132").Value = 4719.706
$ws.Range("I132").Value = 4794.2256
$ws.Range("J132").Value = 3949.6667
$ws.Range("K132").Value = 14382.6768
$ws.Range("L132").Value = 11849.0001
$ws.Range("M132").Value = -11852.6768
$ws.Range("N132").Value = -16909.0001
$ws.Range("H136").Value = 1284.2632
$ws.Range("I136").Value = 1194.2667
$ws.Range("J136").Value = 1621.75
$ws.Range("K136").Value = 3582.800099999999
$ws.Range("L136").Value = 4865.25
$ws.Range("M136").Value = -1032.800099999999
$ws.Range("N136").Value = -9965.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 4000.5
$ws.Range("J31").Value = 4500.6665
$ws.Range("L31").Value = 4500.6665
$ws.Range("N31").Value = -5196.6665
$ws.Range("H51").Value = 11872.238
$ws.Range("I51").Value = 13461.134
$ws.Range("J51").Value = 7900
$ws.Range("K51").Value = 13461.134
$ws.Range("L51").Value = 7900
$ws.Range("M51").Value = -12951.134
$ws.Range("N51").Value = -8920
$ws.Range("H52").Value = 13673.833
$ws.Range("I52").Value = 11999.6
$ws.Range("K52").Value = 11999.6
$ws.Range("M52").Value = -11773.6
$ws.Range("H81").Value = 10658.25
$ws.Range("I81").Value = 16999.75
$ws.Range("J81").Value = 7487.5
$ws.Range("K81").Value = 33999.5
$ws.Range("L81").Value = 14975
$ws.Range("M81").Value = -32938.5
$ws.Range("N81").Value = -17097
$ws.Range("H84").Value = 10658.25
$ws.Range("I84").Value = 16999.75
$ws.Range("J84").Value = 7487.5
$ws.Range("K84").Value = 169997.5
$ws.Range("L84").Value = 74875
$ws.Range("M84").Value = -164693.5
$ws.Range("N84").Value = -85483
$ws.Range("H100").Value = 1013.9286
$ws.Range("I100").Value = 1023.9
$ws.Range("K100").Value = 2047.8
$ws.Range("M100").Value = -1506.8
$ws.Range("H107").Value = 1107.3334
$ws.Range("I107").Value = 999
$ws.Range("J107").Value = 1324
$ws.Range("K107").Value = 2997
$ws.Range("L107").Value = 3972
$ws.Range("M107").Value = -1077
$ws.Range("N107").Value = -7812
$ws.Range("H112").Value = 66961.336
$ws.Range("J112").Value = 66961.336
$ws.Range("L112").Value = 66961.336
$ws.Range("N112").Value = -69915.336
$ws.Range("H122").Value = 9151.023999999999
$ws.Range("I122").Value = 2405.7334
$ws.Range("J122").Value = 27547.273
$ws.Range("K122").Value = 7217.2002
$ws.Range("L122").Value = 82641.819
$ws.Range("M122").Value = -4767.2002
$ws.Range("N122").Value = -87541.819
$ws.Range("H126").Value = 3006.4
$ws.Range("I126").Value = 2502.6365
$ws.Range("J126").Value = 3622.111
$ws.Range("K126").Value = 7507.9095
$ws.Range("L126").Value = 10866.333
$ws.Range("M126").Value = -5037.9095
$ws.Range("N126").Value = -15806.333
$ws.Range("H132").Value = 1347.9565
$ws.Range("I132").Value = 1285.15
$ws.Range("K132").Value = 3855.45
$ws.Range("M132").Value = -1325.45
$ws.Range("H136").Value = 351712.12
$ws.Range("I136").Value = 417811.34
$ws.Range("J136").Value = 87315.336
$ws.Range("K136").Value = 1253434.02
$ws.Range("L136").Value = 261946.008
$ws.Range("M136").Value = -1250884.02
$ws.Range("N136").Value = -267046.008

Write-Host "Applied 557 cell changes"